$d = $word.ActiveDocument

# --- Update the date heading paragraph ---
$dateRange = $d.Paragraphs(1).Range
$dateRange.Text = "2023-08-01 Tuesday"

# --- Update each table cell value (by position, in document order) ---
$t = $d.Tables(1)

$cellRange = $t.Cell(1, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "15+48=63"
$cellRange = $t.Cell(1, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "44+37=81"
$cellRange = $t.Cell(1, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "9+4=13"
$cellRange = $t.Cell(1, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "19+45=64"
$cellRange = $t.Cell(1, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "50-32=18"

$cellRange = $t.Cell(2, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "50-37=13"
$cellRange = $t.Cell(2, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "70-69=1"
$cellRange = $t.Cell(2, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "13+68=81"
$cellRange = $t.Cell(2, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "80-64=16"
$cellRange = $t.Cell(2, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "98-69=29"

$cellRange = $t.Cell(3, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "9+24=33"
$cellRange = $t.Cell(3, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "25+48=73"
$cellRange = $t.Cell(3, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "28+7=35"
$cellRange = $t.Cell(3, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "23+38=61"
$cellRange = $t.Cell(3, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "70-63=7"

$cellRange = $t.Cell(4, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "53+8=61"
$cellRange = $t.Cell(4, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "72-56=16"
$cellRange = $t.Cell(4, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27+66=93"
$cellRange = $t.Cell(4, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "73-67=6"
$cellRange = $t.Cell(4, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "40-18=22"

$cellRange = $t.Cell(5, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "71-17=54"
$cellRange = $t.Cell(5, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "58+8=66"
$cellRange = $t.Cell(5, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "6+29=35"
$cellRange = $t.Cell(5, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "47+29=76"
$cellRange = $t.Cell(5, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "80-9=71"

$cellRange = $t.Cell(6, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "74-68=6"
$cellRange = $t.Cell(6, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "90-18=72"
$cellRange = $t.Cell(6, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "24+49=73"
$cellRange = $t.Cell(6, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "91-89=2"
$cellRange = $t.Cell(6, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "8+17=25"

$cellRange = $t.Cell(7, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "51-47=4"
$cellRange = $t.Cell(7, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "35+36=71"
$cellRange = $t.Cell(7, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "77+19=96"
$cellRange = $t.Cell(7, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "40-26=14"
$cellRange = $t.Cell(7, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "8+14=22"

$cellRange = $t.Cell(8, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "13+79=92"
$cellRange = $t.Cell(8, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "30-18=12"
$cellRange = $t.Cell(8, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "67-9=58"
$cellRange = $t.Cell(8, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "8+13=21"
$cellRange = $t.Cell(8, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27+67=94"

$cellRange = $t.Cell(9, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "17+26=43"
$cellRange = $t.Cell(9, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "16+77=93"
$cellRange = $t.Cell(9, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "97-39=58"
$cellRange = $t.Cell(9, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "71-22=49"
$cellRange = $t.Cell(9, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "79+8=87"

$cellRange = $t.Cell(10, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "83-39=44"
$cellRange = $t.Cell(10, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "37+54=91"
$cellRange = $t.Cell(10, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "70-46=24"
$cellRange = $t.Cell(10, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "8+24=32"
$cellRange = $t.Cell(10, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "82-55=27"

$cellRange = $t.Cell(11, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "26+49=75"
$cellRange = $t.Cell(11, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "50-11=39"
$cellRange = $t.Cell(11, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "39+14=53"
$cellRange = $t.Cell(11, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "6+55=61"
$cellRange = $t.Cell(11, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "92-5=87"

$cellRange = $t.Cell(12, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "54-48=6"
$cellRange = $t.Cell(12, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "47+49=96"
$cellRange = $t.Cell(12, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "45+17=62"
$cellRange = $t.Cell(12, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "92-17=75"
$cellRange = $t.Cell(12, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "28+18=46"

$cellRange = $t.Cell(13, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "48+6=54"
$cellRange = $t.Cell(13, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "78+16=94"
$cellRange = $t.Cell(13, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "66+5=71"
$cellRange = $t.Cell(13, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "55-38=17"
$cellRange = $t.Cell(13, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "38+34=72"

$cellRange = $t.Cell(14, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "82-69=13"
$cellRange = $t.Cell(14, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27+36=63"
$cellRange = $t.Cell(14, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "46-29=17"
$cellRange = $t.Cell(14, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "19+68=87"
$cellRange = $t.Cell(14, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "7+24=31"

$cellRange = $t.Cell(15, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "61-12=49"
$cellRange = $t.Cell(15, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49+2=51"
$cellRange = $t.Cell(15, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "25+39=64"
$cellRange = $t.Cell(15, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "50-33=17"
$cellRange = $t.Cell(15, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "84+8=92"

$cellRange = $t.Cell(16, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "79+18=97"
$cellRange = $t.Cell(16, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "88-29=59"
$cellRange = $t.Cell(16, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "64-35=29"
$cellRange = $t.Cell(16, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "61-17=44"
$cellRange = $t.Cell(16, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "81-19=62"

$cellRange = $t.Cell(17, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "46+25=71"
$cellRange = $t.Cell(17, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "16+49=65"
$cellRange = $t.Cell(17, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "68+23=91"
$cellRange = $t.Cell(17, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "82-9=73"
$cellRange = $t.Cell(17, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27+29=56"

$cellRange = $t.Cell(18, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "27-18=9"
$cellRange = $t.Cell(18, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "87+8=95"
$cellRange = $t.Cell(18, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "19+69=88"
$cellRange = $t.Cell(18, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "43+9=52"
$cellRange = $t.Cell(18, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "9+29=38"

$cellRange = $t.Cell(19, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "30-19=11"
$cellRange = $t.Cell(19, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "7+64=71"
$cellRange = $t.Cell(19, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "71-29=42"
$cellRange = $t.Cell(19, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "26-7=19"
$cellRange = $t.Cell(19, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "98-49=49"

$cellRange = $t.Cell(20, 1).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "49+39=88"
$cellRange = $t.Cell(20, 2).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "88+7=95"
$cellRange = $t.Cell(20, 3).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "7+29=36"
$cellRange = $t.Cell(20, 4).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "29+62=91"
$cellRange = $t.Cell(20, 5).Range
$cellRange.MoveEnd(1, -1) | Out-Null
$cellRange.Text = "23-17=6"

Write-Output "replacements complete"